# Auto-generated Excel COM-interop script
# Applies the 'scheduled runner' data refresh described in the commit diff:
#  - strips the bold/border/center header styling from row 1 of every Leve sheet
#  - rewrites the pricing/profit columns (H:N) for the specific Leve rows whose
#    market data changed, including adding/removing cells where a column's value
#    became populated/blank

$wb = $excel.ActiveWorkbook

# --- 1) Strip header-row styling (bold font / thin border / centered alignment) ---
# on every worksheet's A1:N1 header range, reverting cells to the default style.
foreach ($sheetIndex in 1..$wb.Worksheets.Count) {
    $headerRange = $wb.Worksheets.Item($sheetIndex).Range("A1:N1")
    $headerRange.ClearFormats()
}

# --- 2) Update per-sheet market-price / profit data ---

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(43, 8).Value = 614.2857  # H43: 575 -> 614.2857
$ws.Cells.Item(43, 10).Value = 614.2857  # J43: 575 -> 614.2857
$ws.Cells.Item(43, 12).Value = 614.2857  # L43: 575 -> 614.2857
$ws.Cells.Item(43, 14).Value = -752.2857  # N43: -713 -> -752.2857
$ws.Cells.Item(51, 8).Value = 6120  # H51: 4793.75 -> 6120
$ws.Cells.Item(51, 9).Value = 6740  # I51: 6450 -> 6740
$ws.Cells.Item(51, 10).Value = 5500  # J51: 3800 -> 5500
$ws.Cells.Item(51, 11).Value = 6740  # K51: 6450 -> 6740
$ws.Cells.Item(51, 12).Value = 5500  # L51: 3800 -> 5500
$ws.Cells.Item(51, 13).Value = -6256  # M51: -5966 -> -6256
$ws.Cells.Item(51, 14).Value = -6468  # N51: -4768 -> -6468
$ws.Cells.Item(69, 8).Value = 3650  # H69: 1761.6666 -> 3650
$ws.Cells.Item(69, 10).Value = 0  # J69: 1590 -> 0
$ws.Cells.Item(69, 12).Value = 0  # L69: 4770 -> 0
$ws.Cells.Item(69, 14).ClearContents()  # N69: -6518 -> (removed)
$ws.Cells.Item(72, 8).Value = 3650  # H72: 1761.6666 -> 3650
$ws.Cells.Item(72, 10).Value = 0  # J72: 1590 -> 0
$ws.Cells.Item(72, 12).Value = 0  # L72: 14310 -> 0
$ws.Cells.Item(72, 14).ClearContents()  # N72: -23046 -> (removed)
$ws.Cells.Item(80, 8).Value = 3502770.8  # H80: 8086.1113 -> 3502770.8
$ws.Cells.Item(80, 9).Value = 2583.3333  # I80: 1750 -> 2583.3333
$ws.Cells.Item(80, 10).Value = 4055432  # J80: 10523.077 -> 4055432
$ws.Cells.Item(80, 11).Value = 7749.999899999999  # K80: 5250 -> 7749.999899999999
$ws.Cells.Item(80, 12).Value = 12166296  # L80: 31569.231 -> 12166296
$ws.Cells.Item(80, 13).Value = -6751.999899999999  # M80: -4252 -> -6751.999899999999
$ws.Cells.Item(80, 14).Value = -12168292  # N80: -33565.231 -> -12168292
$ws.Cells.Item(83, 8).Value = 3502770.8  # H83: 8086.1113 -> 3502770.8
$ws.Cells.Item(83, 9).Value = 2583.3333  # I83: 1750 -> 2583.3333
$ws.Cells.Item(83, 10).Value = 4055432  # J83: 10523.077 -> 4055432
$ws.Cells.Item(83, 11).Value = 23249.9997  # K83: 15750 -> 23249.9997
$ws.Cells.Item(83, 12).Value = 36498888  # L83: 94707.693 -> 36498888
$ws.Cells.Item(83, 13).Value = -18257.9997  # M83: -10758 -> -18257.9997
$ws.Cells.Item(83, 14).Value = -36508872  # N83: -104691.693 -> -36508872
$ws.Cells.Item(127, 8).Value = 2162.6  # H127: 1512.3334 -> 2162.6
$ws.Cells.Item(127, 9).Value = 0  # I127: 898.5 -> 0
$ws.Cells.Item(127, 10).Value = 2162.6  # J127: 1687.7142 -> 2162.6
$ws.Cells.Item(127, 11).Value = 0  # K127: 2695.5 -> 0
$ws.Cells.Item(127, 12).Value = 6487.799999999999  # L127: 5063.142599999999 -> 6487.799999999999
$ws.Cells.Item(127, 13).ClearContents()  # M127: 2264.5 -> (removed)
$ws.Cells.Item(127, 14).Value = -16407.8  # N127: -14983.1426 -> -16407.8
$ws.Cells.Item(129, 8).Value = 698.44446  # H129: 814.3333 -> 698.44446
$ws.Cells.Item(129, 9).Value = 547.6667  # I129: 557.8 -> 547.6667
$ws.Cells.Item(129, 10).Value = 1000  # J129: 894.5 -> 1000
$ws.Cells.Item(129, 11).Value = 1643.0001  # K129: 1673.4 -> 1643.0001
$ws.Cells.Item(129, 12).Value = 3000  # L129: 2683.5 -> 3000
$ws.Cells.Item(129, 13).Value = 3356.9999  # M129: 3326.6 -> 3356.9999
$ws.Cells.Item(129, 14).Value = -13000  # N129: -12683.5 -> -13000
$ws.Cells.Item(138, 8).Value = 2323.8965  # H138: 2434.3125 -> 2323.8965
$ws.Cells.Item(138, 9).Value = 647.0714  # I138: 688.9231 -> 647.0714
$ws.Cells.Item(138, 10).Value = 3888.9333  # J138: 3628.5264 -> 3888.9333
$ws.Cells.Item(138, 11).Value = 1941.2142  # K138: 2066.7693 -> 1941.2142
$ws.Cells.Item(138, 12).Value = 11666.7999  # L138: 10885.5792 -> 11666.7999
$ws.Cells.Item(138, 13).Value = 3198.7858  # M138: 3073.2307 -> 3198.7858
$ws.Cells.Item(138, 14).Value = -21946.7999  # N138: -21165.5792 -> -21946.7999
$ws.Cells.Item(141, 8).Value = 2333.3333  # H141: 1700.7894 -> 2333.3333
$ws.Cells.Item(141, 9).Value = 2333.3333  # I141: 1077.3529 -> 2333.3333
$ws.Cells.Item(141, 10).Value = 0  # J141: 7000 -> 0
$ws.Cells.Item(141, 11).Value = 6999.999899999999  # K141: 3232.0587 -> 6999.999899999999
$ws.Cells.Item(141, 12).Value = 0  # L141: 21000 -> 0
$ws.Cells.Item(141, 13).Value = -1819.999899999999  # M141: 1947.9413 -> -1819.999899999999
$ws.Cells.Item(141, 14).ClearContents()  # N141: -31360 -> (removed)

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 3065.1064  # H32: 3170.2 -> 3065.1064
$ws.Cells.Item(32, 9).Value = 1958.4857  # I32: 2034.697 -> 1958.4857
$ws.Cells.Item(32, 10).Value = 6292.75  # J32: 6292.8335 -> 6292.75
$ws.Cells.Item(32, 11).Value = 1958.4857  # K32: 2034.697 -> 1958.4857
$ws.Cells.Item(32, 12).Value = 6292.75  # L32: 6292.8335 -> 6292.75
$ws.Cells.Item(32, 13).Value = -1671.4857  # M32: -1747.697 -> -1671.4857
$ws.Cells.Item(32, 14).Value = -6866.75  # N32: -6866.8335 -> -6866.75
$ws.Cells.Item(46, 8).Value = 8430  # H46: 0 -> 8430
$ws.Cells.Item(46, 9).Value = 8138  # I46: 0 -> 8138
$ws.Cells.Item(46, 10).Value = 8576  # J46: 0 -> 8576
$ws.Cells.Item(46, 11).Value = 8138  # K46: 0 -> 8138
$ws.Cells.Item(46, 12).Value = 8576  # L46: 0 -> 8576
$ws.Cells.Item(46, 13).Value = -7819  # M46: None -> -7819
$ws.Cells.Item(46, 14).Value = -9214  # N46: None -> -9214
$ws.Cells.Item(61, 8).Value = 3138.9333  # H61: 3292 -> 3138.9333
$ws.Cells.Item(61, 9).Value = 2552.4546  # I61: 2777.7 -> 2552.4546
$ws.Cells.Item(61, 10).Value = 4751.75  # J61: 4434.8887 -> 4751.75
$ws.Cells.Item(61, 11).Value = 2552.4546  # K61: 2777.7 -> 2552.4546
$ws.Cells.Item(61, 12).Value = 4751.75  # L61: 4434.8887 -> 4751.75
$ws.Cells.Item(61, 13).Value = -2340.4546  # M61: -2565.7 -> -2340.4546
$ws.Cells.Item(61, 14).Value = -5175.75  # N61: -4858.8887 -> -5175.75
$ws.Cells.Item(122, 8).Value = 3009.625  # H122: 2241.0625 -> 3009.625
$ws.Cells.Item(122, 9).Value = 3533.3333  # I122: 1968 -> 3533.3333
$ws.Cells.Item(122, 10).Value = 2695.4  # J122: 2696.1667 -> 2695.4
$ws.Cells.Item(122, 11).Value = 10599.9999  # K122: 5904 -> 10599.9999
$ws.Cells.Item(122, 12).Value = 8086.200000000001  # L122: 8088.500100000001 -> 8086.200000000001
$ws.Cells.Item(122, 13).Value = -8149.999899999999  # M122: -3454 -> -8149.999899999999
$ws.Cells.Item(122, 14).Value = -12986.2  # N122: -12988.5001 -> -12986.2
$ws.Cells.Item(136, 8).Value = 3138.9333  # H136: 3292 -> 3138.9333
$ws.Cells.Item(136, 9).Value = 2552.4546  # I136: 2777.7 -> 2552.4546
$ws.Cells.Item(136, 10).Value = 4751.75  # J136: 4434.8887 -> 4751.75
$ws.Cells.Item(136, 11).Value = 7657.3638  # K136: 8333.099999999999 -> 7657.3638
$ws.Cells.Item(136, 12).Value = 14255.25  # L136: 13304.6661 -> 14255.25
$ws.Cells.Item(136, 13).Value = -5107.3638  # M136: -5783.099999999999 -> -5107.3638
$ws.Cells.Item(136, 14).Value = -19355.25  # N136: -18404.6661 -> -19355.25

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 1250  # H20: 980 -> 1250
$ws.Cells.Item(20, 9).Value = 1100  # I20: 800 -> 1100
$ws.Cells.Item(20, 10).Value = 1400  # J20: 1700 -> 1400
$ws.Cells.Item(20, 11).Value = 1100  # K20: 800 -> 1100
$ws.Cells.Item(20, 12).Value = 1400  # L20: 1700 -> 1400
$ws.Cells.Item(20, 13).Value = -853  # M20: -553 -> -853
$ws.Cells.Item(20, 14).Value = -1894  # N20: -2194 -> -1894
$ws.Cells.Item(81, 8).Value = 14974.714  # H81: 15985.5 -> 14974.714
$ws.Cells.Item(81, 10).Value = 14974.714  # J81: 15985.5 -> 14974.714
$ws.Cells.Item(81, 12).Value = 14974.714  # L81: 15985.5 -> 14974.714
$ws.Cells.Item(81, 14).Value = -17096.714  # N81: -18107.5 -> -17096.714
$ws.Cells.Item(84, 8).Value = 14974.714  # H84: 15985.5 -> 14974.714
$ws.Cells.Item(84, 10).Value = 14974.714  # J84: 15985.5 -> 14974.714
$ws.Cells.Item(84, 12).Value = 44924.142  # L84: 47956.5 -> 44924.142
$ws.Cells.Item(84, 14).Value = -55532.142  # N84: -58564.5 -> -55532.142
$ws.Cells.Item(105, 8).Value = 2390  # H105: 2220 -> 2390
$ws.Cells.Item(105, 9).Value = 3150  # I105: 2300 -> 3150
$ws.Cells.Item(105, 11).Value = 3150  # K105: 2300 -> 3150
$ws.Cells.Item(105, 13).Value = -1403  # M105: -553 -> -1403
$ws.Cells.Item(134, 8).Value = 2892.5625  # H134: 2253.0444 -> 2892.5625
$ws.Cells.Item(134, 9).Value = 2922  # I134: 2272.5476 -> 2922
$ws.Cells.Item(134, 11).Value = 8766  # K134: 6817.6428 -> 8766
$ws.Cells.Item(134, 13).Value = -6231  # M134: -4282.6428 -> -6231

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 1224.5  # H16: 1399.5 -> 1224.5
$ws.Cells.Item(16, 9).Value = 1066  # I16: 1099.5 -> 1066
$ws.Cells.Item(16, 10).Value = 1700  # J16: 1699.5 -> 1700
$ws.Cells.Item(16, 11).Value = 1066  # K16: 1099.5 -> 1066
$ws.Cells.Item(16, 12).Value = 1700  # L16: 1699.5 -> 1700
$ws.Cells.Item(16, 13).Value = -779  # M16: -812.5 -> -779
$ws.Cells.Item(16, 14).Value = -2274  # N16: -2273.5 -> -2274
$ws.Cells.Item(22, 8).Value = 634.44446  # H22: 587.1 -> 634.44446
$ws.Cells.Item(22, 9).Value = 535  # I22: 481.57144 -> 535
$ws.Cells.Item(22, 11).Value = 535  # K22: 481.57144 -> 535
$ws.Cells.Item(22, 13).Value = -185  # M22: -131.57144 -> -185
$ws.Cells.Item(99, 8).Value = 4849  # H99: 4755.04 -> 4849
$ws.Cells.Item(99, 10).Value = 7020  # J99: 6568 -> 7020
$ws.Cells.Item(99, 12).Value = 7020  # L99: 6568 -> 7020
$ws.Cells.Item(99, 14).Value = -10016  # N99: -9564 -> -10016
$ws.Cells.Item(113, 8).Value = 1224.5  # H113: 1399.5 -> 1224.5
$ws.Cells.Item(113, 9).Value = 1066  # I113: 1099.5 -> 1066
$ws.Cells.Item(113, 10).Value = 1700  # J113: 1699.5 -> 1700
$ws.Cells.Item(113, 11).Value = 1066  # K113: 1099.5 -> 1066
$ws.Cells.Item(113, 12).Value = 1700  # L113: 1699.5 -> 1700
$ws.Cells.Item(113, 13).Value = 1104  # M113: 1070.5 -> 1104
$ws.Cells.Item(113, 14).Value = -6040  # N113: -6039.5 -> -6040
$ws.Cells.Item(126, 8).Value = 4849  # H126: 4755.04 -> 4849
$ws.Cells.Item(126, 10).Value = 7020  # J126: 6568 -> 7020
$ws.Cells.Item(126, 12).Value = 21060  # L126: 19704 -> 21060
$ws.Cells.Item(126, 14).Value = -26000  # N126: -24644 -> -26000

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(58, 8).Value = 3251  # H58: 3276.75 -> 3251
$ws.Cells.Item(58, 10).Value = 3963.6  # J58: 4004.8 -> 3963.6
$ws.Cells.Item(58, 12).Value = 11890.8  # L58: 12014.4 -> 11890.8
$ws.Cells.Item(58, 14).Value = -12146.8  # N58: -12270.4 -> -12146.8
$ws.Cells.Item(68, 8).Value = 1573.1818  # H68: 1475 -> 1573.1818
$ws.Cells.Item(68, 10).Value = 1573.1818  # J68: 1475 -> 1573.1818
$ws.Cells.Item(68, 12).Value = 4719.5454  # L68: 4425 -> 4719.5454
$ws.Cells.Item(68, 14).Value = -6341.5454  # N68: -6047 -> -6341.5454
$ws.Cells.Item(71, 8).Value = 1573.1818  # H71: 1475 -> 1573.1818
$ws.Cells.Item(71, 10).Value = 1573.1818  # J71: 1475 -> 1573.1818
$ws.Cells.Item(71, 12).Value = 14158.6362  # L71: 13275 -> 14158.6362
$ws.Cells.Item(71, 14).Value = -22270.6362  # N71: -21387 -> -22270.6362
$ws.Cells.Item(131, 8).Value = 787.75  # H131: 792.56 -> 787.75
$ws.Cells.Item(131, 9).Value = 396.44446  # I131: 398.75 -> 396.44446
$ws.Cells.Item(131, 10).Value = 826.45056  # J131: 826.8043 -> 826.45056
$ws.Cells.Item(131, 11).Value = 1189.33338  # K131: 1196.25 -> 1189.33338
$ws.Cells.Item(131, 12).Value = 2479.35168  # L131: 2480.4129 -> 2479.35168
$ws.Cells.Item(131, 13).Value = 3850.66662  # M131: 3843.75 -> 3850.66662
$ws.Cells.Item(131, 14).Value = -12559.35168  # N131: -12560.4129 -> -12559.35168

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(97, 8).Value = 1721.4231  # H97: 1780.68 -> 1721.4231
$ws.Cells.Item(97, 9).Value = 963.6667  # I97: 1006.2353 -> 963.6667
$ws.Cells.Item(97, 11).Value = 963.6667  # K97: 1006.2353 -> 963.6667
$ws.Cells.Item(97, 13).Value = -467.6667  # M97: -510.2353000000001 -> -467.6667
$ws.Cells.Item(126, 8).Value = 4566  # H126: 4809.304 -> 4566
$ws.Cells.Item(126, 9).Value = 3692.8572  # I126: 4008.3333 -> 3692.8572
$ws.Cells.Item(126, 10).Value = 5677.273  # J126: 5683.091 -> 5677.273
$ws.Cells.Item(126, 11).Value = 11078.5716  # K126: 12024.9999 -> 11078.5716
$ws.Cells.Item(126, 12).Value = 17031.819  # L126: 17049.273 -> 17031.819
$ws.Cells.Item(126, 13).Value = -8608.5716  # M126: -9554.999899999999 -> -8608.5716
$ws.Cells.Item(126, 14).Value = -21971.819  # N126: -21989.273 -> -21971.819

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 3275  # H81: 1775.375 -> 3275
$ws.Cells.Item(81, 9).Value = 2366.6667  # I81: 1325.5 -> 2366.6667
$ws.Cells.Item(81, 10).Value = 6000  # J81: 3125 -> 6000
$ws.Cells.Item(81, 11).Value = 4733.3334  # K81: 2651 -> 4733.3334
$ws.Cells.Item(81, 12).Value = 12000  # L81: 6250 -> 12000
$ws.Cells.Item(81, 13).Value = -3672.3334  # M81: -1590 -> -3672.3334
$ws.Cells.Item(81, 14).Value = -14122  # N81: -8372 -> -14122
$ws.Cells.Item(84, 8).Value = 3275  # H84: 1775.375 -> 3275
$ws.Cells.Item(84, 9).Value = 2366.6667  # I84: 1325.5 -> 2366.6667
$ws.Cells.Item(84, 10).Value = 6000  # J84: 3125 -> 6000
$ws.Cells.Item(84, 11).Value = 23666.667  # K84: 13255 -> 23666.667
$ws.Cells.Item(84, 12).Value = 60000  # L84: 31250 -> 60000
$ws.Cells.Item(84, 13).Value = -18362.667  # M84: -7951 -> -18362.667
$ws.Cells.Item(84, 14).Value = -70608  # N84: -41858 -> -70608
$ws.Cells.Item(122, 8).Value = 1800  # H122: 2097.3076 -> 1800
$ws.Cells.Item(122, 9).Value = 1800  # I122: 1898.75 -> 1800
$ws.Cells.Item(122, 10).Value = 0  # J122: 2415 -> 0
$ws.Cells.Item(122, 11).Value = 5400  # K122: 5696.25 -> 5400
$ws.Cells.Item(122, 12).Value = 0  # L122: 7245 -> 0
$ws.Cells.Item(122, 13).Value = -2950  # M122: -3246.25 -> -2950
$ws.Cells.Item(122, 14).ClearContents()  # N122: -12145 -> (removed)
$ws.Cells.Item(132, 8).Value = 1621.4722  # H132: 1716.7354 -> 1621.4722
$ws.Cells.Item(132, 9).Value = 1400.625  # I132: 1527.7727 -> 1400.625
$ws.Cells.Item(132, 11).Value = 4201.875  # K132: 4583.3181 -> 4201.875
$ws.Cells.Item(132, 13).Value = -1671.875  # M132: -2053.3181 -> -1671.875
